# TRIMAZKON address list update
# - ip_address_list: flip several "favorite" flags to TRUE/1, convert two
#   flag cells to boolean type, fix a typo in the 503_Witte note and
#   re-insert a duplicated 497_Edcha row that had been merged into 503_Witte.
# - ip_address_fav_list: populate with the rows that are now flagged as
#   favorites on ip_address_list.
# - disk_list: drop the now-unneeded hyperlink/style on C1 and reorder two
#   blocks of rows (518_Valeo/515_ZF now come before 514_Teleflex/474_B Austin).
# - projects_bin2 (hidden): populate with a small snapshot of rows pulled
#   from the other two sheets.

$wb = $excel.ActiveWorkbook

$wsIp   = $wb.Worksheets.Item("ip_address_list")
$wsFav  = $wb.Worksheets.Item("ip_address_fav_list")
$wsDisk = $wb.Worksheets.Item("disk_list")
$wsBin  = $wb.Worksheets.Item("projects_bin2")

# --- capture "before" values we will need later, before anything moves ---
$ip_row1 = $wsIp.Range("A1:E1").Value2   # Domaci Wifi
$ip_row3 = $wsIp.Range("A3:E3").Value2   # 529_Witte
$ip_row4 = $wsIp.Range("A4:E4").Value2   # 527_Teijin
$ip_row5 = $wsIp.Range("A5:E5").Value2   # 518_Valeo II
$ip_row9 = $wsIp.Range("A9:E9").Value2   # 514_Teleflex

$disk_row6 = $wsDisk.Range("A6:E6").Value2   # 515_ZF (A:E only, F is blank anyway)

# ================= ip_address_list =================

# Flip plain favourite flags 0 -> 1
$wsIp.Range("E1").Value2 = 1   # Domaci Wifi
$wsIp.Range("E3").Value2 = 1   # 529_Witte
$wsIp.Range("E4").Value2 = 1   # 527_Teijin
$wsIp.Range("E5").Value2 = 1   # 518_Valeo II
$wsIp.Range("E9").Value2 = 1   # 514_Teleflex

# Two flags become real booleans (still FALSE/0)
$wsIp.Range("E7").Value2 = $false   # 515_ZF Stara kkkBoleslav
$wsIp.Range("E8").Value2 = $false   # 515_

# Row 11 (497_Edcha) / row 14 (503_Witte) get untangled: the 503_Witte note
# (with a typo fixed) moves up to row 11, and a fresh row is inserted at 12
# to hold the 497_Edcha entry that used to live there, shifting the old
# rows 12/13 down to 13/14. The old duplicate row 14 is then removed.
$row11old = $wsIp.Range("A11:E11").Value2   # 497_Edcha ...
$row14old = $wsIp.Range("A14:E14").Value2   # 503_Witte ... (typo'd note)

$wsIp.Rows.Item(12).Insert()

$wsIp.Range("A12:E12").Value2 = $row11old
$wsIp.Range("A11:E11").Value2 = $row14old

# fix the "175k" typo introduced when the note is restored into row 11
$fixed = $wsIp.Range("D11").Value2 -replace "175k", "175"
$wsIp.Range("D11").Value2 = $fixed

# the old row 14 got pushed down to row 15 by the insert above; drop it
$wsIp.Rows.Item(15).Delete()

# ================= ip_address_fav_list =================
# New favourites list mirrors the rows just flagged above (in the order
# they appear on ip_address_list), each with its flag forced to 1.

$wsFav.Range("A1:E1").Value2 = $ip_row9    # 514_Teleflex
$wsFav.Range("E1").Value2 = 1

$wsFav.Range("A2:E2").Value2 = $ip_row5    # 518_Valeo II
$wsFav.Range("E2").Value2 = 1

$wsFav.Range("A3:E3").Value2 = $ip_row4    # 527_Teijin
$wsFav.Range("E3").Value2 = 1

$wsFav.Range("A4:E4").Value2 = $ip_row3    # 529_Witte
$wsFav.Range("E4").Value2 = 1

$wsFav.Range("A5:E5").Value2 = $ip_row1    # Domaci Wifi
$wsFav.Range("E5").Value2 = 1

# ================= disk_list =================

# C1 no longer carries an external hyperlink / hyperlink styling
$wsDisk.Range("C1").Hyperlinks.Delete()
$wsDisk.Range("C1").Style = "Normal"

# Swap the (514_Teleflex, 474_B Austin) block with the (515_ZF, 518_Valeo)
# block, each keeping its own internal order, so 518_Valeo/515_ZF now come
# first, followed by 514_Teleflex/474_B Austin.
$old4 = $wsDisk.Range("A4:F4").Value2
$old5 = $wsDisk.Range("A5:F5").Value2
$old6 = $wsDisk.Range("A6:F6").Value2
$old7 = $wsDisk.Range("A7:F7").Value2

$wsDisk.Range("A4:F4").Value2 = $old7
$wsDisk.Range("A5:F5").Value2 = $old6
$wsDisk.Range("A6:F6").Value2 = $old4
$wsDisk.Range("A7:F7").Value2 = $old5

# ================= projects_bin2 (hidden) =================

$wsBin.Range("A1:E1").Value2 = $ip_row3     # 529_Witte, still flagged 0 here
$wsBin.Range("A2:E2").Value2 = $disk_row6   # 515_ZF
$wsBin.Range("A3:E3").Value2 = $ip_row1     # Domaci Wifi
$wsBin.Range("E3").Value2 = 1
